$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B from 56 to 57 (stored width = ColumnWidth + 5/6, so back-solve)
$ws.Columns.Item(2).ColumnWidth = 56.166666666666664

# Update row 2 values
$ws.Range("A2").Value = "#1429"
$ws.Range("B2").Value = "A Product Support Excel T-shirt 2D lần 15 0310 - Black"
$ws.Range("E2").Value = "05:21 03/10/2024"
$ws.Range("H2").Value = "luyenstore_6450966757448"
# Leading apostrophe keeps this purely-numeric-looking value stored as text
# (matches the source cell's Text type instead of being auto-coerced to a number)
$ws.Range("I2").Value = "'40416861683784"
$ws.Range("P2").Value = "https://cdn.shopify.com/s/files/1/0013/7602/6696/files/c9d00636ce2b818b96e07fa5087aa300_3f44b37b-66a1-4375-94c8-21357320a95e.jpg?v=1727950833"
